$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name3"
$ws.Range("C1").Value = "Address3"
$ws.Range("D1").Value = "Website3"
$ws.Range("F1").Value = "Latitude3"
$ws.Range("G1").Value = "Longitude3"
$ws.Range("H1").Value = "Description3"
$ws.Range("I1").Value = "Hours3"

$ws.Range("C2").Select()
